$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 3875.96
$ws.Range("J40").Value = 4599.8335
$ws.Range("L40").Value = 4599.8335
$ws.Range("N40").Value = -4949.8335
# Row 132
$ws.Range("H132").Value = 2886.3022
$ws.Range("I132").Value = 2941.3076
$ws.Range("K132").Value = 8823.9228
$ws.Range("M132").Value = -6293.9228
# Row 141
$ws.Range("H141").Value = 7632.1875
$ws.Range("J141").Value = 3047
$ws.Range("L141").Value = 9141
$ws.Range("N141").Value = -19501

$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 190.75
$ws.Range("I5").Value = 190.75
$ws.Range("K5").Value = 190.75
$ws.Range("M5").Value = -78.75
# Row 32
$ws.Range("H32").Value = 5217.57
$ws.Range("I32").Value = 4281.8105
$ws.Range("J32").Value = 22997
$ws.Range("K32").Value = 4281.8105
$ws.Range("L32").Value = 22997
$ws.Range("M32").Value = -3994.8105
$ws.Range("N32").Value = -23571
# Row 63
$ws.Range("H63").Value = 5796.706
$ws.Range("I63").Value = 3571.7778
$ws.Range("J63").Value = 8299.75
$ws.Range("K63").Value = 3571.7778
$ws.Range("L63").Value = 8299.75
$ws.Range("M63").Value = -2885.7778
$ws.Range("N63").Value = -9671.75
# Row 66
$ws.Range("H66").Value = 5796.706
$ws.Range("I66").Value = 3571.7778
$ws.Range("J66").Value = 8299.75
$ws.Range("K66").Value = 17858.889
$ws.Range("L66").Value = 41498.75
$ws.Range("M66").Value = -14426.889
$ws.Range("N66").Value = -48362.75
# Row 132
$ws.Range("H132").Value = 25722.5
$ws.Range("I132").Value = 1744.0741
$ws.Range("J132").Value = 84578.63
$ws.Range("K132").Value = 5232.2223
$ws.Range("L132").Value = 253735.89
$ws.Range("M132").Value = -2702.2223
$ws.Range("N132").Value = -258795.89
# Row 134
$ws.Range("H134").Value = 64898.5
$ws.Range("J134").Value = 64898.5
$ws.Range("L134").Value = 64898.5
$ws.Range("N134").Value = -75038.5

$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 190.75
$ws.Range("I4").Value = 190.75
$ws.Range("K4").Value = 190.75
$ws.Range("M4").Value = -75.75
# Row 86
$ws.Range("H86").Value = 12749145
$ws.Range("I86").Value = 21669480
$ws.Range("J86").Value = 5808.7144
$ws.Range("K86").Value = 21669480
$ws.Range("L86").Value = 5808.7144
$ws.Range("M86").Value = -21668357
$ws.Range("N86").Value = -8054.7144
# Row 89
$ws.Range("H89").Value = 12749145
$ws.Range("I89").Value = 21669480
$ws.Range("J89").Value = 5808.7144
$ws.Range("K89").Value = 108347400
$ws.Range("L89").Value = 29043.572
$ws.Range("M89").Value = -108341784
$ws.Range("N89").Value = -40275.572

$ws = $wb.Worksheets.Item("CRP")
# Row 23
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
# Row 26
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").ClearContents()
# Row 27
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
# Row 31
$ws.Range("H31").Value = 17835.955
$ws.Range("I31").Value = 2699.0232
$ws.Range("K31").Value = 2699.0232
$ws.Range("M31").Value = -2404.0232
# Row 34
$ws.Range("H34").Value = 17835.955
$ws.Range("I34").Value = 2699.0232
$ws.Range("K34").Value = 2699.0232
$ws.Range("M34").Value = -2497.0232
# Row 36
$ws.Range("H36").Value = 10000
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
# Row 40
$ws.Range("H40").Value = 10000
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
# Row 99
$ws.Range("H99").Value = 3406.3447
$ws.Range("I99").Value = 3119.652
$ws.Range("K99").Value = 3119.652
$ws.Range("M99").Value = -1621.652
# Row 105
$ws.Range("H105").Value = 816.36365
$ws.Range("I105").Value = 795.2222
$ws.Range("J105").Value = 911.5
$ws.Range("K105").Value = 795.2222
$ws.Range("L105").Value = 911.5
$ws.Range("M105").Value = 951.7778
$ws.Range("N105").Value = -4405.5
# Row 126
$ws.Range("H126").Value = 3406.3447
$ws.Range("I126").Value = 3119.652
$ws.Range("K126").Value = 9358.956
$ws.Range("M126").Value = -6888.956
# Row 132
$ws.Range("H132").Value = 47264.117
$ws.Range("I132").Value = 36555.4
$ws.Range("K132").Value = 109666.2
$ws.Range("M132").Value = -107136.2
# Row 141
$ws.Range("H141").Value = 236667.33
$ws.Range("J141").Value = 236667.33
$ws.Range("L141").Value = 236667.33
$ws.Range("N141").Value = -247027.33

$ws = $wb.Worksheets.Item("GSM")
# Row 18
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
# Row 43
$ws.Range("H43").Value = 12239.625
$ws.Range("I43").Value = 7986.1665
$ws.Range("K43").Value = 7986.1665
$ws.Range("M43").Value = -7835.1665
# Row 46
$ws.Range("H46").Value = 13000
$ws.Range("J46").Value = 25000
$ws.Range("L46").Value = 25000
$ws.Range("N46").Value = -25312
# Row 122
$ws.Range("H122").Value = 242773.75
$ws.Range("I122").Value = 271667
$ws.Range("K122").Value = 815001
$ws.Range("M122").Value = -812551
# Row 132
$ws.Range("H132").Value = 3079.4492
$ws.Range("I132").Value = 2774.074
$ws.Range("J132").Value = 4178.8
$ws.Range("K132").Value = 8322.222
$ws.Range("L132").Value = 12536.4
$ws.Range("M132").Value = -5792.222
$ws.Range("N132").Value = -17596.4

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 75773.164
$ws.Range("I22").Value = 81752.55
$ws.Range("K22").Value = 81752.55
$ws.Range("M22").Value = -81457.55
# Row 27
$ws.Range("H27").Value = 75773.164
$ws.Range("I27").Value = 81752.55
$ws.Range("K27").Value = 81752.55
$ws.Range("M27").Value = -81645.55
# Row 55
$ws.Range("H55").Value = 2055.2942
$ws.Range("I55").Value = 1994.8334
$ws.Range("K55").Value = 1994.8334
$ws.Range("M55").Value = -1821.8334
# Row 136
$ws.Range("H136").Value = 72262.31
$ws.Range("I136").Value = 90152.74000000001
$ws.Range("J136").Value = 3682.3333
$ws.Range("K136").Value = 270458.22
$ws.Range("L136").Value = 11046.9999
$ws.Range("M136").Value = -267908.22
$ws.Range("N136").Value = -16146.9999

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 8267.549000000001
$ws.Range("J62").Value = 8722.111000000001
$ws.Range("L62").Value = 8722.111000000001
$ws.Range("N62").Value = -9970.111000000001
# Row 65
$ws.Range("H65").Value = 8267.549000000001
$ws.Range("J65").Value = 8722.111000000001
$ws.Range("L65").Value = 43610.55500000001
$ws.Range("N65").Value = -49850.55500000001
# Row 69
$ws.Range("H69").Value = 20000
$ws.Range("J69").Value = 20000
$ws.Range("L69").Value = 20000
$ws.Range("N69").Value = -21498
# Row 72
$ws.Range("H72").Value = 20000
$ws.Range("J72").Value = 20000
$ws.Range("L72").Value = 60000
$ws.Range("N72").Value = -67488
# Row 81
$ws.Range("H81").Value = 10423260
$ws.Range("I81").Value = 13893856
$ws.Range("K81").Value = 27787712
$ws.Range("M81").Value = -27786651
# Row 84
$ws.Range("H84").Value = 10423260
$ws.Range("I84").Value = 13893856
$ws.Range("K84").Value = 138938560
$ws.Range("M84").Value = -138933256

